# The spec's "gather" command section (8 bullet points, from
# `Implement a "gather" command` through the `Syntax: ...` paragraph
# ending in "...Armor or Weapon") gets marked done the same way the
# earlier, already-finished bullets are: colored green (RGB 00B050),
# matching e.g. "Implement a command to create a Mine location" /
# "...Forest location" right above it.

$d = $word.ActiveDocument

# wdColor packs RGB 00B050 as 0x00BBGGRR.
$green = 5287936

$paragraphs = $d.Paragraphs
$total = $paragraphs.Count

$startIndex = -1
for ($i = 1; $i -le $total; $i++) {
    $text = $paragraphs.Item($i).Range.Text
    if ($text -like "Implement a*gather*command*") {
        $startIndex = $i
        break
    }
}

if ($startIndex -gt 0) {
    for ($i = $startIndex; $i -le $total; $i++) {
        $text = $paragraphs.Item($i).Range.Text
        if ($i -gt $startIndex -and $text -like "Implement a*craft*command*") {
            break
        }
        $paragraphs.Item($i).Range.Font.Color = $green
    }
}
